$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume/1h (E) updates reflecting the latest market snapshot
$ws.Range("D2").Value = "'27.766.11"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "'1.851.28"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  -1.56%  "
$ws.Range("D5").Value = "'318.31"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").Value = "'1.010"
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("D7").Value = "'0.4309"
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("D8").Value = "'0.3754"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "'0.07349"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").Value = "'0.8765"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "'1.862.99"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "'6.747"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "'5.442"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "'0.07124"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "'89.04"
$ws.Range("E16").Value = "  +4.66%  "
$ws.Range("D17").Value = "'1.013"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("D18").Value = "'0.000009012"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "'1.010"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").Value = "'15.44"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").Value = "'27.785.13"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "'5.218"
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("D23").Value = "'11.05"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").Value = "'2.077.62"
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("D25").Value = "'1.973"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").Value = "'155.35"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").Value = "'2.163"
$ws.Range("E28").Value = "  +9.18%  "
$ws.Range("D29").Value = "'5.360"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").Value = "'118.83"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").Value = "'0.08933"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").Value = "'0.7781"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").Value = "'4.545"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'2.923"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("D36").Value = "'1.011"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("D37").Value = "'1.133"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").Value = "'0.01983"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("D39").Value = "'0.05346"
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("D40").Value = "'2.897"
$ws.Range("E40").Value = "  +1.95%  "
$ws.Range("D41").Value = "'7.161"
$ws.Range("E41").Value = "  +4.70%  "
$ws.Range("E42").Value = "  +1.37%  "
$ws.Range("D43").Value = "'0.5139"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").Value = "'8.805"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "'10.74"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").Value = "'107.59"
$ws.Range("E46").Value = "  -2.15%  "
$ws.Range("D47").Value = "'0.4765"
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("D48").Value = "'0.06473"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D51").Value = "'1.853"
$ws.Range("E51").Value = "  -2.36%  "

# Rows 49-50 swap rank order: NEARProtocol now ranks above PaxDollar
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.692"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'1.011"
$ws.Range("E50").Value = "  -1.46%  "
